$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (16) down into the new row (17)
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row's values
$ws.Range("A17").Value = "train.csv"
$ws.Range("B17").Value = "Kaggle"
$ws.Range("C17").Value = "https://www.kaggle.com/c/rossmann-store-sales/data"

# Turn the link cell into a real hyperlink
$ws.Hyperlinks.Add($ws.Range("C17"), "https://www.kaggle.com/c/rossmann-store-sales/data")

# Hyperlinks.Add re-stamps the cell style; restore the look-and-feel of the row above
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)  # xlPasteFormats

# Match row height used by the other single-line rows
$ws.Rows(17).RowHeight = $ws.Rows(7).RowHeight

# Reset the view: no frozen/scrolled top-left cell, and select C7:C10 like the target
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C7:C10").Select()
